$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cells: _old -> _FV2410, _new -> _FV2504
$oldSuffix = "_old"
$newSuffix = "_new"

for ($col = 1; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $val = $cell.Value
    if ($val -ne $null) {
        if ($val.ToString().EndsWith($oldSuffix)) {
            $base = $val.ToString().Substring(0, $val.ToString().Length - $oldSuffix.Length)
            $cell.Value = "$base" + "_FV2410"
        } elseif ($val.ToString().EndsWith($newSuffix)) {
            $base = $val.ToString().Substring(0, $val.ToString().Length - $newSuffix.Length)
            $cell.Value = "$base" + "_FV2504"
        }
    }
}
